$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link / volume(%) text updates
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("E3").Value = '  +1.20%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("E12").Value = '  -1.68%  '
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("E21").Value = '  -1.35%  '
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("E27").Value = '  +1.17%  '
$ws.Range("E28").Value = '  -4.14%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E32").Value = '  +4.86%  '
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("E34").Value = '  -2.71%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("E38").Value = '  -3.17%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E39").Value = '  +2.39%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E40").Value = '  +11.11%  '
$ws.Range("E41").Value = '  -2.80%  '
$ws.Range("E42").Value = '  -3.19%  '
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("E44").Value = '  -1.89%  '
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("E46").Value = '  -4.03%  '
$ws.Range("E47").Value = '  -5.17%  '
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("E51").Value = '  -1.84%  '

# Price column (D) updates: set as Text first so values like "1.00" / "0.890" /
# "66.903.35" keep their exact literal formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("D2").Value = '66.903.35'
$ws.Range("D3").Value = '3.518.33'
$ws.Range("D5").Value = '585.58'
$ws.Range("D6").Value = '177.34'
$ws.Range("D8").Value = '0.602'
$ws.Range("D9").Value = '3.516.03'
$ws.Range("D11").Value = '6.92'
$ws.Range("D12").Value = '0.425'
$ws.Range("D13").Value = '4.120.89'
$ws.Range("D14").Value = '30.66'
$ws.Range("D16").Value = '66.862.05'
$ws.Range("D17").Value = '0.0000175'
$ws.Range("D18").Value = '3.515.78'
$ws.Range("D20").Value = '14.07'
$ws.Range("D21").Value = '381.38'
$ws.Range("D25").Value = '0.538'
$ws.Range("D26").Value = '71.52'
$ws.Range("D29").Value = '0.173'
$ws.Range("D31").Value = '6.02'
$ws.Range("D32").Value = '24.61'
$ws.Range("D35").Value = '1.00'
$ws.Range("D37").Value = '1.57'
$ws.Range("D38").Value = '158.83'
$ws.Range("D39").Value = '0.890'
$ws.Range("D40").Value = '29.02'
$ws.Range("D41").Value = '1.81'
$ws.Range("D42").Value = '2.63'
$ws.Range("D44").Value = '4.54'
$ws.Range("D46").Value = '2.717.72'
$ws.Range("D47").Value = '25.69'
$ws.Range("D48").Value = '40.60'
$ws.Range("D50").Value = '328.05'
